$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45179 = 2023-09-10) for every
# data row (2-176). Bump it by one day to 45180 (2023-09-11) across the whole range.
$ws.Range("C2:C176").Value = 45180
